# ------------------------------------------------------------------
# OptionsNotFound.xlsx - add two new BT option rows (3900C / 4000C)
# derived from the same underlying (EMINI S&P DEC0), inserted above
# the existing BABA/VIX option rows.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 3 (pushes the existing rows 3-4 down to 5-6,
# carrying their formatting/styles with them; the new rows 3-4 inherit the
# style of the row above, same as Excel's native Insert behaviour).
$ws.Rows("3:4").Insert()

# Populate the columns that introduce brand-new instrument identifiers first,
# entered per-column across both new rows (as they would be pasted/derived
# together from the source ticker lookup), so shared-string ordering matches.
$ws.Cells.Item(3,5).Value2 = "O_9995ee58-d18a-11ea-a100-402102c1e5a3_8ES DEC0 3900 C | 18 DEC 20 | C | A | 3900.0000000000_TBRX"
$ws.Cells.Item(4,5).Value2 = "O_999600aa-d18a-11ea-a100-402102c1e5a3_8ES DEC0 4000 C | 18 DEC 20 | C | A | 4000.0000000000_TBRX"
$ws.Cells.Item(3,21).Value2 = "8ES DEC0 3900 C | 18 DEC 20 | C | A | 3900.0000000000"
$ws.Cells.Item(4,21).Value2 = "8ES DEC0 4000 C | 18 DEC 20 | C | A | 4000.0000000000"
$ws.Cells.Item(3,42).Value2 = "8ES DEC0 3900 C"
$ws.Cells.Item(4,42).Value2 = "8ES DEC0 4000 C"

# Fill in the remaining fields for each new row.
# Row 3
$ws.Cells.Item(3,1).Value2 = "USD"
$ws.Cells.Item(3,2).Value2 = "Options"
$ws.Cells.Item(3,3).Value2 = "Call"
$ws.Cells.Item(3,4).Value2 = "BT"
$ws.Cells.Item(3,6).Value2 = -500
$ws.Cells.Item(3,7).Value2 = -134516.64853989449
$ws.Cells.Item(3,8).Value2 = -11.5360928557
$ws.Cells.Item(3,9).Value2 = -39356.272427341799
$ws.Cells.Item(3,10).Value2 = 13739.524550309599
$ws.Cells.Item(3,11).Value2 = -21192.1063684484
$ws.Cells.Item(3,12).Value2 = "A"
$ws.Cells.Item(3,14).Value2 = "A"
$ws.Cells.Item(3,15).Value2 = 44183
$ws.Cells.Item(3,16).Value2 = "DEC 20"
$ws.Cells.Item(3,17).Value2 = 0
$ws.Cells.Item(3,18).Value2 = "C"
$ws.Cells.Item(3,19).Value2 = "P"
$ws.Cells.Item(3,22).Value2 = "N"
$ws.Cells.Item(3,23).Value2 = 50
$ws.Cells.Item(3,24).Value2 = 44074
$ws.Cells.Item(3,25).Value2 = "TBRX"
$ws.Cells.Item(3,26).Value2 = 30200831
$ws.Cells.Item(3,27).Value2 = 3900
$ws.Cells.Item(3,28).Value2 = "L_S_9baccdce-3e91-11ea-b3f0-05b9e229bfa0_EMINI S&P DEC0 | DEC 20_TBRX"
$ws.Cells.Item(3,29).Value2 = "EMINI S&P DEC0 | DEC 20"
$ws.Cells.Item(3,37).Value2 = "FI_OPTION"
$ws.Cells.Item(3,38).Value2 = "Y"
$ws.Cells.Item(3,45).Value2 = "Equities"
$ws.Cells.Item(3,46).Value2 = "Equities"
$ws.Cells.Item(3,47).Value2 = "Options"
$ws.Cells.Item(3,48).Value2 = "Call"
$ws.Cells.Item(3,50).Value2 = 1
$ws.Cells.Item(3,51).Value2 = 0
$ws.Cells.Item(3,53).Value2 = "Call"
$ws.Cells.Item(3,54).Value2 = "Equities"
$ws.Cells.Item(3,56).Value2 = "Options"

# Row 4
$ws.Cells.Item(4,1).Value2 = "USD"
$ws.Cells.Item(4,2).Value2 = "Options"
$ws.Cells.Item(4,3).Value2 = "Call"
$ws.Cells.Item(4,4).Value2 = "BT"
$ws.Cells.Item(4,6).Value2 = 225
$ws.Cells.Item(4,7).Value2 = 42478.193191086997
$ws.Cells.Item(4,8).Value2 = 3.9768287984000001
$ws.Cells.Item(4,9).Value2 = 12460.797135027
$ws.Cells.Item(4,10).Value2 = -4927.1101480056004
$ws.Cells.Item(4,11).Value2 = 7613.0491956186997
$ws.Cells.Item(4,12).Value2 = "A"
$ws.Cells.Item(4,14).Value2 = "A"
$ws.Cells.Item(4,15).Value2 = 44183
$ws.Cells.Item(4,16).Value2 = "DEC 20"
$ws.Cells.Item(4,17).Value2 = 0
$ws.Cells.Item(4,18).Value2 = "C"
$ws.Cells.Item(4,19).Value2 = "P"
$ws.Cells.Item(4,22).Value2 = "N"
$ws.Cells.Item(4,23).Value2 = 50
$ws.Cells.Item(4,24).Value2 = 44074
$ws.Cells.Item(4,25).Value2 = "TBRX"
$ws.Cells.Item(4,26).Value2 = 30200831
$ws.Cells.Item(4,27).Value2 = 4000
$ws.Cells.Item(4,28).Value2 = "L_S_9baccdce-3e91-11ea-b3f0-05b9e229bfa0_EMINI S&P DEC0 | DEC 20_TBRX"
$ws.Cells.Item(4,29).Value2 = "EMINI S&P DEC0 | DEC 20"
$ws.Cells.Item(4,37).Value2 = "FI_OPTION"
$ws.Cells.Item(4,38).Value2 = "Y"
$ws.Cells.Item(4,45).Value2 = "Equities"
$ws.Cells.Item(4,46).Value2 = "Equities"
$ws.Cells.Item(4,47).Value2 = "Options"
$ws.Cells.Item(4,48).Value2 = "Call"
$ws.Cells.Item(4,50).Value2 = 1
$ws.Cells.Item(4,51).Value2 = 0
$ws.Cells.Item(4,53).Value2 = "Call"
$ws.Cells.Item(4,54).Value2 = "Equities"
$ws.Cells.Item(4,56).Value2 = "Options"

# Leave the sheet scrolled/selected where the user ended up working: further
# right (around the ticker/underlying columns) and with the full new block
# of rows selected.
$ws.Range("A7:XFD21").Select()
